# Updates the cryptos price/volume table (sheet1) to the latest scraped
# snapshot: refreshed "Price" (D) and "Volume(1h)" (E) figures for every
# coin row, plus a couple of rows whose rank swapped places (Maker <->
# InternetComputer(DFINITY) at rows 33/34, RenderToken <-> FraxShare at
# rows 44/45), which moves each coin's Name/Link/Price/Volume down or up
# one row.
$changes = @(
    @{ Row=2; Col="D"; Value='27.958.21' }
    @{ Row=2; Col="E"; Value='  +0.22%  ' }
    @{ Row=3; Col="D"; Value='1.639.54' }
    @{ Row=3; Col="E"; Value='  +0.04%  ' }
    @{ Row=4; Col="E"; Value='  -0.02%  ' }
    @{ Row=5; Col="D"; Value='212.59' }
    @{ Row=5; Col="E"; Value='  +0.05%  ' }
    @{ Row=6; Col="E"; Value='  -0.21%  ' }
    @{ Row=7; Col="E"; Value='  -0.02%  ' }
    @{ Row=8; Col="D"; Value='23.41' }
    @{ Row=8; Col="E"; Value='  +0.04%  ' }
    @{ Row=9; Col="D"; Value='0.259' }
    @{ Row=9; Col="E"; Value='  -2.33%  ' }
    @{ Row=10; Col="E"; Value='  +0.17%  ' }
    @{ Row=11; Col="E"; Value='  +1.86%  ' }
    @{ Row=12; Col="D"; Value='1.872.18' }
    @{ Row=12; Col="E"; Value='  +0.04%  ' }
    @{ Row=13; Col="D"; Value='1.638.17' }
    @{ Row=13; Col="E"; Value='  -0.09%  ' }
    @{ Row=14; Col="E"; Value='  +0.32%  ' }
    @{ Row=15; Col="E"; Value='  +1.35%  ' }
    @{ Row=16; Col="E"; Value='  -0.14%  ' }
    @{ Row=17; Col="D"; Value='27.955.13' }
    @{ Row=17; Col="E"; Value='  +0.26%  ' }
    @{ Row=18; Col="D"; Value='233.17' }
    @{ Row=18; Col="E"; Value='  +0.68%  ' }
    @{ Row=19; Col="E"; Value='  -0.19%  ' }
    @{ Row=20; Col="D"; Value='7.56' }
    @{ Row=20; Col="E"; Value='  -1.07%  ' }
    @{ Row=21; Col="E"; Value='  -0.02%  ' }
    @{ Row=22; Col="D"; Value='10.45' }
    @{ Row=22; Col="E"; Value='  -2.48%  ' }
    @{ Row=23; Col="E"; Value='  -0.06%  ' }
    @{ Row=24; Col="E"; Value='  -2.97%  ' }
    @{ Row=25; Col="D"; Value='152.94' }
    @{ Row=25; Col="E"; Value='  +1.38%  ' }
    @{ Row=26; Col="D"; Value='6.92' }
    @{ Row=26; Col="E"; Value='  +0.29%  ' }
    @{ Row=27; Col="D"; Value='15.66' }
    @{ Row=27; Col="E"; Value='  -0.22%  ' }
    @{ Row=28; Col="E"; Value='  -0.58%  ' }
    @{ Row=29; Col="E"; Value='  -0.09%  ' }
    @{ Row=30; Col="E"; Value='  +0.14%  ' }
    @{ Row=31; Col="E"; Value='  +0.45%  ' }
    @{ Row=32; Col="D"; Value='3.40' }
    @{ Row=32; Col="E"; Value='  +2.98%  ' }
    @{ Row=33; Col="B"; Value='InternetComputer(DFINITY)' }
    @{ Row=33; Col="C"; Value='https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp' }
    @{ Row=33; Col="D"; Value='3.09' }
    @{ Row=33; Col="E"; Value='  +0.03%  ' }
    @{ Row=34; Col="B"; Value='Maker' }
    @{ Row=34; Col="C"; Value='https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr' }
    @{ Row=34; Col="D"; Value='1.404.72' }
    @{ Row=34; Col="E"; Value='  -3.65%  ' }
    @{ Row=35; Col="E"; Value='  +1.88%  ' }
    @{ Row=36; Col="E"; Value='  +1.16%  ' }
    @{ Row=37; Col="E"; Value='  +0.65%  ' }
    @{ Row=38; Col="E"; Value='  +0.19%  ' }
    @{ Row=39; Col="D"; Value='0.928' }
    @{ Row=39; Col="E"; Value='  +0.19%  ' }
    @{ Row=40; Col="E"; Value='  -0.91%  ' }
    @{ Row=41; Col="E"; Value='  +0.68%  ' }
    @{ Row=42; Col="E"; Value='  -0.12%  ' }
    @{ Row=43; Col="E"; Value='  -3.03%  ' }
    @{ Row=44; Col="B"; Value='FraxShare' }
    @{ Row=44; Col="C"; Value='https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs' }
    @{ Row=44; Col="D"; Value='5.52' }
    @{ Row=44; Col="E"; Value='  +3.29%  ' }
    @{ Row=45; Col="B"; Value='RenderToken' }
    @{ Row=45; Col="C"; Value='https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr' }
    @{ Row=45; Col="D"; Value='1.84' }
    @{ Row=45; Col="E"; Value='  +1.81%  ' }
    @{ Row=46; Col="E"; Value='  -0.19%  ' }
    @{ Row=47; Col="D"; Value='1.781.02' }
    @{ Row=47; Col="E"; Value='  -0.49%  ' }
    @{ Row=48; Col="D"; Value='87.95' }
    @{ Row=48; Col="E"; Value='  -0.09%  ' }
    @{ Row=49; Col="E"; Value='  -0.07%  ' }
    @{ Row=50; Col="D"; Value='0.0505' }
    @{ Row=50; Col="E"; Value='  -0.23%  ' }
    @{ Row=51; Col="D"; Value='7.59' }
    @{ Row=51; Col="E"; Value='  -1.94%  ' }
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($chg in $changes) {
    $cell = $ws.Cells.Item($chg.Row, [int][char]$chg.Col - [int][char]'A' + 1)
    if ($chg.Col -eq "D") {
        # Column D holds price strings that can look numeric (e.g. "3.40", "0.0505").
        # Force them to stay plain text (matching the source inlineStr cells) by using
        # the quote-prefix convention, then restore the cell's default style.
        $cell.Value = "'" + $chg.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $chg.Value
    }
}
